$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.379.24'
$ws.Range("E2").Value = '  +3.06%  '
$ws.Range("D3").Value = '1.791.89'
$ws.Range("E3").Value = '  +3.83%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.34%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '335.95'
$ws.Range("E5").Value = '  +1.09%  '
$ws.Range("E6").Value = '  +0.43%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3781'
$ws.Range("E7").Value = '  +1.91%  '
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3431'
$ws.Range("E8").Value = '  +2.07%  '
$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.54'
$ws.Range("E9").Value = '  +0.06%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.201'
$ws.Range("E10").Value = '  +1.40%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07489'
$ws.Range("E11").Value = '  +1.18%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.003'
$ws.Range("E12").Value = '  +0.29%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.96'
$ws.Range("E13").Value = '  +9.20%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.481'
$ws.Range("E14").Value = '  +1.50%  '
$ws.Range("D15").Value = '1.790.64'
$ws.Range("E15").Value = '  +3.78%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.034'
$ws.Range("E16").Value = '  -0.34%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001093'
$ws.Range("E17").Value = '  +2.09%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06618'
$ws.Range("E18").Value = '  -0.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '84.53'
$ws.Range("E19").Value = '  +3.22%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  +0.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.35'
$ws.Range("E21").Value = '  +4.87%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.463'
$ws.Range("E22").Value = '  +5.38%  '
$ws.Range("D23").Value = '27.376.67'
$ws.Range("E23").Value = '  +3.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.56'
$ws.Range("E24").Value = '  -1.65%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.460'
$ws.Range("E25").Value = '  +0.95%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.568'
$ws.Range("E26").Value = '  +7.10%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.41'
$ws.Range("E27").Value = '  +10.15%  '
$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.478'
$ws.Range("E28").Value = '  +3.79%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '150.34'
$ws.Range("E29").Value = '  -0.91%  '
$ws.Range("D30").Value = '1.996.03'
$ws.Range("E30").Value = '  +4.11%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '133.58'
$ws.Range("E31").Value = '  +1.58%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.052'
$ws.Range("E32").Value = '  -1.59%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.123'
$ws.Range("E33").Value = '  +2.46%  '
$ws.Range("E34").Value = '  +0.93%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '13.28'
$ws.Range("E35").Value = '  +4.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.669'
$ws.Range("E36").Value = '  -1.28%  '
$ws.Range("B37").Value = 'TheSandbox'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6885'
$ws.Range("E37").Value = '  +11.19%  '
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.431'
$ws.Range("E38").Value = '  +1.36%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06353'
$ws.Range("E39").Value = '  +2.34%  '
$ws.Range("E40").Value = '  +2.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.810'
$ws.Range("E41").Value = '  +4.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.02341'
$ws.Range("E42").Value = '  +0.64%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.268'
$ws.Range("E43").Value = '  +4.12%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.38'
$ws.Range("E44").Value = '  +0.60%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6451'
$ws.Range("E45").Value = '  +7.24%  '
$ws.Range("E46").Value = '  +0.40%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.850'
$ws.Range("E47").Value = '  -1.44%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.117'
$ws.Range("E48").Value = '  +3.71%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '129.75'
$ws.Range("E49").Value = '  +0.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07214'
$ws.Range("E50").Value = '  +0.51%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.45'
$ws.Range("E51").Value = '  +3.34%  '
